$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting of the neighboring header
# cell (G1) so it picks up the same bold/border/alignment style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding numeric value in H2 (no special formatting, like B2:G2)
$ws.Range("H2").Value = 0
